# Delete row 796 entirely, shifting all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(796).Delete()
